$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 9836.125
$ws.Range("I28").Value = 5489
$ws.Range("J28").Value = 19399.8
$ws.Range("K28").Value = 5489
$ws.Range("L28").Value = 19399.8
$ws.Range("M28").Value = -5004
$ws.Range("N28").Value = -20369.8
$ws.Range("H31").Value = 700.75
$ws.Range("J31").Value = 2000
$ws.Range("L31").Value = 6000
$ws.Range("N31").Value = -6460
$ws.Range("H64").Value = 6000.8
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 6000.8
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 6000.8
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -6496.8
$ws.Range("H67").Value = 6000.8
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 6000.8
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 6000.8
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -7716.8
$ws.Range("H74").Value = 9136.182000000001
$ws.Range("I74").Value = 8099.6
$ws.Range("K74").Value = 8099.6
$ws.Range("M74").Value = -7163.6
$ws.Range("H77").Value = 9136.182000000001
$ws.Range("I77").Value = 8099.6
$ws.Range("K77").Value = 40498
$ws.Range("M77").Value = -35818
$ws.Range("H92").Value = 1649.9333
$ws.Range("I92").Value = 1742.36
$ws.Range("J92").Value = 1187.8
$ws.Range("K92").Value = 1742.36
$ws.Range("L92").Value = 1187.8
$ws.Range("M92").Value = -494.3599999999999
$ws.Range("N92").Value = -3683.8
$ws.Range("H101").Value = 2523.6
$ws.Range("I101").Value = 2496.6667
$ws.Range("J101").Value = 2564
$ws.Range("K101").Value = 7490.000100000001
$ws.Range("L101").Value = 7692
$ws.Range("M101").Value = -5868.000100000001
$ws.Range("N101").Value = -10936
$ws.Range("H113").Value = 3979.7
$ws.Range("I113").Value = 3849.25
$ws.Range("K113").Value = 3849.25
$ws.Range("M113").Value = -595.25
$ws.Range("H132").Value = 1482.7142
$ws.Range("I132").Value = 1243.1666
$ws.Range("J132").Value = 1913.9
$ws.Range("K132").Value = 3729.4998
$ws.Range("L132").Value = 5741.700000000001
$ws.Range("M132").Value = -1199.4998
$ws.Range("N132").Value = -10801.7
$ws.Range("H138").Value = 3091
$ws.Range("I138").Value = 2290.5
$ws.Range("K138").Value = 6871.5
$ws.Range("M138").Value = -1731.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4842.231
$ws.Range("I61").Value = 3546.5789
$ws.Range("K61").Value = 3546.5789
$ws.Range("M61").Value = -3334.5789
$ws.Range("H63").Value = 5333.3335
$ws.Range("I63").Value = 5500
$ws.Range("K63").Value = 5500
$ws.Range("M63").Value = -4814
$ws.Range("H66").Value = 5333.3335
$ws.Range("I66").Value = 5500
$ws.Range("K66").Value = 27500
$ws.Range("M66").Value = -24068
$ws.Range("H122").Value = 2133.3333
$ws.Range("I122").Value = 2133.3333
$ws.Range("K122").Value = 6399.999899999999
$ws.Range("M122").Value = -3949.999899999999
$ws.Range("H136").Value = 4842.231
$ws.Range("I136").Value = 3546.5789
$ws.Range("K136").Value = 10639.7367
$ws.Range("M136").Value = -8089.736699999999

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H86").Value = 127486.69
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 127486.69
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 127486.69
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -129732.69
$ws.Range("H89").Value = 127486.69
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 127486.69
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 637433.45
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -648665.45
$ws.Range("H94").Value = 1575.381
$ws.Range("I94").Value = 1240.1765
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 1240.1765
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -789.1765
$ws.Range("N94").Value = -3902
$ws.Range("H105").Value = 62517296
$ws.Range("I105").Value = 66684950
$ws.Range("K105").Value = 66684950
$ws.Range("M105").Value = -66683203

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3960.4348
$ws.Range("I16").Value = 1499
$ws.Range("K16").Value = 1499
$ws.Range("M16").Value = -1212
$ws.Range("H58").Value = 3544.8572
$ws.Range("J58").Value = 7049.875
$ws.Range("L58").Value = 7049.875
$ws.Range("N58").Value = -7455.875
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 3960.4348
$ws.Range("I113").Value = 1499
$ws.Range("K113").Value = 1499
$ws.Range("M113").Value = 671
$ws.Range("H122").Value = 1050.6923
$ws.Range("I122").Value = 1059.909
$ws.Range("K122").Value = 3179.727
$ws.Range("M122").Value = -729.7270000000003
$ws.Range("H134").Value = 8409.556
$ws.Range("I134").Value = 8148.1904
$ws.Range("K134").Value = 24444.5712
$ws.Range("M134").Value = -21909.5712
$ws.Range("H136").Value = 3544.8572
$ws.Range("J136").Value = 7049.875
$ws.Range("L136").Value = 21149.625
$ws.Range("N136").Value = -26249.625

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 9236.308000000001
$ws.Range("I64").Value = 9996.333000000001
$ws.Range("K64").Value = 29988.999
$ws.Range("M64").Value = -29718.999
$ws.Range("H67").Value = 9236.308000000001
$ws.Range("I67").Value = 9996.333000000001
$ws.Range("K67").Value = 29988.999
$ws.Range("M67").Value = -29052.999

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 685.8125
$ws.Range("I97").Value = 824
$ws.Range("J97").Value = 381.8
$ws.Range("K97").Value = 824
$ws.Range("L97").Value = 381.8
$ws.Range("M97").Value = -328
$ws.Range("N97").Value = -1373.8
$ws.Range("H107").Value = 271.14285
$ws.Range("I107").Value = 279.8
$ws.Range("J107").Value = 249.5
$ws.Range("K107").Value = 279.8
$ws.Range("L107").Value = 249.5
$ws.Range("M107").Value = 1640.2
$ws.Range("N107").Value = -4089.5
$ws.Range("H113").Value = 38450.125
$ws.Range("I113").Value = 1902.5
$ws.Range("J113").Value = 74997.75
$ws.Range("K113").Value = 1902.5
$ws.Range("L113").Value = 74997.75
$ws.Range("M113").Value = 267.5
$ws.Range("N113").Value = -79337.75

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2805.068
$ws.Range("I22").Value = 1678.0741
$ws.Range("J22").Value = 4595
$ws.Range("K22").Value = 1678.0741
$ws.Range("L22").Value = 4595
$ws.Range("M22").Value = -1383.0741
$ws.Range("N22").Value = -5185
$ws.Range("H27").Value = 2805.068
$ws.Range("I27").Value = 1678.0741
$ws.Range("J27").Value = 4595
$ws.Range("K27").Value = 1678.0741
$ws.Range("L27").Value = 4595
$ws.Range("M27").Value = -1571.0741
$ws.Range("N27").Value = -4809
$ws.Range("H40").Value = 13001.3125
$ws.Range("I40").Value = 11850.615
$ws.Range("K40").Value = 11850.615
$ws.Range("M40").Value = -11714.615
$ws.Range("H69").Value = 125000
$ws.Range("J69").Value = 125000
$ws.Range("L69").Value = 125000
$ws.Range("N69").Value = -126622
$ws.Range("H72").Value = 125000
$ws.Range("J72").Value = 125000
$ws.Range("L72").Value = 375000
$ws.Range("N72").Value = -383112
$ws.Range("H82").Value = 2372.5833
$ws.Range("I82").Value = 1855.8
$ws.Range("J82").Value = 2741.7144
$ws.Range("K82").Value = 1855.8
$ws.Range("L82").Value = 2741.7144
$ws.Range("M82").Value = -1494.8
$ws.Range("N82").Value = -3463.7144
$ws.Range("H85").Value = 2372.5833
$ws.Range("I85").Value = 1855.8
$ws.Range("J85").Value = 2741.7144
$ws.Range("K85").Value = 1855.8
$ws.Range("L85").Value = 2741.7144
$ws.Range("M85").Value = -607.8
$ws.Range("N85").Value = -5237.7144
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H132").Value = 4685.7026
$ws.Range("I132").Value = 3803.8215
$ws.Range("K132").Value = 11411.4645
$ws.Range("M132").Value = -8881.4645

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 16166.667
$ws.Range("J5").Value = 16166.667
$ws.Range("L5").Value = 16166.667
$ws.Range("N5").Value = -16390.667
$ws.Range("H62").Value = 47832.668
$ws.Range("I62").Value = 66747.5
$ws.Range("J62").Value = 10003
$ws.Range("K62").Value = 66747.5
$ws.Range("L62").Value = 10003
$ws.Range("M62").Value = -66123.5
$ws.Range("N62").Value = -11251
$ws.Range("H65").Value = 47832.668
$ws.Range("I65").Value = 66747.5
$ws.Range("J65").Value = 10003
$ws.Range("K65").Value = 333737.5
$ws.Range("L65").Value = 50015
$ws.Range("M65").Value = -330617.5
$ws.Range("N65").Value = -56255
$ws.Range("H103").Value = 73392.336
$ws.Range("J103").Value = 73392.336
$ws.Range("L103").Value = 73392.336
$ws.Range("N103").Value = -75736.336
